# Update "last updated" timestamp banner (A1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 11:33"

# Swap the Caceres / Gran Canaria rows (row 31 / row 32) since Gran Canaria's
# updated total now outranks Caceres in this cases-descending sorted table.
$ws.Range("A31").Value = "Gran Canaria"
$ws.Range("A32").Value = "Caceres"

# Refresh the per-province numbers (Casos totales / Casos activos / Recuperados / Muertes)
$ws.Range("B4").Value  = 62395
$ws.Range("C4").Value  = 37808
$ws.Range("D4").Value  = 16211
$ws.Range("E4").Value  = 8376

$ws.Range("B5").Value  = 50366
$ws.Range("C5").Value  = 21898
$ws.Range("D5").Value  = 23248
$ws.Range("E5").Value  = 5220

$ws.Range("B6").Value  = 17334
$ws.Range("C6").Value  = 6900
$ws.Range("D6").Value  = 8616
$ws.Range("E6").Value  = 1818

$ws.Range("B7").Value  = 16050
$ws.Range("C7").Value  = 5794
$ws.Range("D7").Value  = 7666
$ws.Range("E7").Value  = 2590

$ws.Range("B9").Value  = 12194
$ws.Range("C9").Value  = 7021
$ws.Range("D9").Value  = 3910
$ws.Range("E9").Value  = 1263

$ws.Range("B10").Value = 9011
$ws.Range("C10").Value = 6234
$ws.Range("D10").Value = 2208
$ws.Range("E10").Value = 569

$ws.Range("B13").Value = 5188
$ws.Range("C13").Value = 2624
$ws.Range("D13").Value = 1802
$ws.Range("E13").Value = 762

$ws.Range("B15").Value = 4918
$ws.Range("C15").Value = 2494
$ws.Range("D15").Value = 1958
$ws.Range("E15").Value = 466

$ws.Range("B16").Value = 3961
$ws.Range("C16").Value = 2312
$ws.Range("D16").Value = 1314
$ws.Range("E16").Value = 335

$ws.Range("B23").Value = 2849
$ws.Range("C23").Value = 2090
$ws.Range("D23").Value = 299
$ws.Range("E23").Value = 460

$ws.Range("B30").Value = 2306
$ws.Range("C30").Value = 917
$ws.Range("D30").Value = 1105
$ws.Range("E30").Value = 284

$ws.Range("B31").Value = 2225
$ws.Range("C31").Value = 1166
$ws.Range("D31").Value = 918
$ws.Range("E31").Value = 141

$ws.Range("B32").Value = 2220
$ws.Range("C32").Value = 422
$ws.Range("D32").Value = 1482
$ws.Range("E32").Value = 316

$ws.Range("B33").Value = 2206
$ws.Range("C33").Value = 1610
$ws.Range("D33").Value = 399
$ws.Range("E33").Value = 197

$ws.Range("B38").Value = 1492
$ws.Range("C38").Value = 1277
$ws.Range("D38").Value = 81
$ws.Range("E38").Value = 134

$ws.Range("B59").Value = 119
$ws.Range("C59").Value = 101
$ws.Range("D59").Value = 16
